# Add Open-DNA-Collections support:
#  1. Insert a new "OpenDNACollectionsSource" sheet right before
#     "GenomeCoordinatesSource" (i.e. right after "IGEMSource"), with the
#     same shape as the other single-repository-id source sheets.
#  2. Add "open_dna_collections" to the repository-name validation list on
#     every sheet that already offers the addgene/genbank/.../seva list.

$wb = $excel.ActiveWorkbook

$newList = """addgene,genbank,benchling,snapgene,euroscarf,igem,wekwikgene,seva,open_dna_collections"""

function Set-RepoValidation($ws, $colLetter) {
    $range = $ws.Range($colLetter + "2:" + $colLetter + "1048576")
    $range.Validation.Delete()
    $range.Validation.Add(3, 1, 1, $newList)
    $range.Validation.IgnoreBlank = $true
    $range.Validation.InCellDropdown = $true
    $range.Validation.ShowInput = $false
    $range.Validation.ShowError = $false
}

# --- 1. Create the new sheet, inserted before GenomeCoordinatesSource ---
$beforeSheet = $wb.Worksheets.Item("GenomeCoordinatesSource")
$newSheet = $wb.Worksheets.Add($beforeSheet)
$newSheet.Name = "OpenDNACollectionsSource"

$headers = @("sequence_file_url", "repository_id", "repository_name", "type", "output_name", "database_id", "input", "id")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# repository_name is column C on the new sheet
Set-RepoValidation $newSheet "C"

# --- 2. Extend the existing repository-name validation lists ---
Set-RepoValidation $wb.Worksheets.Item("RepositoryIdSource") "B"
Set-RepoValidation $wb.Worksheets.Item("AddgeneIdSource") "D"
Set-RepoValidation $wb.Worksheets.Item("WekWikGeneIdSource") "C"
Set-RepoValidation $wb.Worksheets.Item("SEVASource") "C"
Set-RepoValidation $wb.Worksheets.Item("BenchlingUrlSource") "B"
Set-RepoValidation $wb.Worksheets.Item("SnapGenePlasmidSource") "B"
Set-RepoValidation $wb.Worksheets.Item("EuroscarfSource") "B"
Set-RepoValidation $wb.Worksheets.Item("IGEMSource") "C"

Write-Host "OpenDNACollectionsSource sheet added and validations updated"
